$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 matching style of existing header cells (copy format from E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "Modelo"

# Update slightly changed numeric values
$ws.Range("B2").Value = 0.395593553237457
$ws.Range("D2").Value = 0.4321333824756294

# Add new model name cell
$ws.Range("F2").Value = "Pipeline(steps=[('model', GradientBoostingRegressor(n_estimators=150))])"
